$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.510.32'
$ws.Range("E2").Value = '  -2.64%  '
$ws.Range("D3").Value = '1.671.32'
$ws.Range("E4").Value = '  +0.39%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.61'
$ws.Range("E5").Value = '  -1.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5154'
$ws.Range("E6").Value = '  -2.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.007'
$ws.Range("E7").Value = '  +0.39%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06480'
$ws.Range("E8").Value = '  -1.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2578'
$ws.Range("E9").Value = '  -2.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.05'
$ws.Range("E10").Value = '  -3.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07673'
$ws.Range("E11").Value = '  +0.69%  '
$ws.Range("D12").Value = '1.675.92'
$ws.Range("E12").Value = '  -1.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.350'
$ws.Range("E13").Value = '  -4.63%  '
$ws.Range("D14").Value = '1.901.16'
$ws.Range("E14").Value = '  -1.83%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5586'
$ws.Range("E15").Value = '  -2.09%  '
$ws.Range("D16").Value = '0.0₅8070'
$ws.Range("E16").Value = '  -0.91%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.91'
$ws.Range("E17").Value = '  -3.75%  '
$ws.Range("D18").Value = '26.542.18'
$ws.Range("E18").Value = '  -2.53%  '
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '211.47'
$ws.Range("E20").Value = '  -1.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.447'
$ws.Range("E21").Value = '  -4.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.14'
$ws.Range("E22").Value = '  -2.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.911'
$ws.Range("E23").Value = '  -0.47%  '
$ws.Range("E24").Value = '  +0.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.54'
$ws.Range("E25").Value = '  +2.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.739'
$ws.Range("E26").Value = '  -0.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1169'
$ws.Range("E27").Value = '  -3.69%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.024'
$ws.Range("E28").Value = '  -2.88%  '
$ws.Range("E29").Value = '  -2.94%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05227'
$ws.Range("E30").Value = '  -2.59%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.264'
$ws.Range("E31").Value = '  -1.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.380'
$ws.Range("E32").Value = '  -3.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.223'
$ws.Range("E33").Value = '  -5.58%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.583'
$ws.Range("E34").Value = '  -3.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.774'
$ws.Range("E35").Value = '  -3.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.376'
$ws.Range("E36").Value = '  -1.73%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9266'
$ws.Range("E37").Value = '  -1.85%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5740'
$ws.Range("E38").Value = '  -1.54%  '
$ws.Range("D39").Value = '1.165.06'
$ws.Range("E39").Value = '  +11.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01601'
$ws.Range("E40").Value = '  -1.57%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8507'
$ws.Range("E41").Value = '  +1.75%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.007'
$ws.Range("E42").Value = '  +0.43%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.641'
$ws.Range("E43").Value = '  -3.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.28'
$ws.Range("E44").Value = '  -0.38%  '
$ws.Range("D45").Value = '1.811.23'
$ws.Range("E45").Value = '  -1.78%  '
$ws.Range("D46").Value = '0.0₈115'
$ws.Range("E46").Value = '  +0.46%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4491'
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '56.04'
$ws.Range("E48").Value = '  -3.10%  '
$ws.Range("E49").Value = '  +0.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.950'
$ws.Range("E50").Value = '  -1.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05130'
$ws.Range("E51").Value = '  -2.05%  '
